$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$rng = $ws.Range("A1:C3")
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(7).Weight = 2
$rng.Borders.Item(8).LineStyle = 1
$rng.Borders.Item(8).Weight = 2
$rng.Borders.Item(9).LineStyle = 1
$rng.Borders.Item(9).Weight = 2
$rng.Borders.Item(10).LineStyle = 1
$rng.Borders.Item(10).Weight = 2
